$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing "Location" column (E) one column to the right (F),
# then drop in the new "Input_group" column in its place (E).
for ($r = 1; $r -le 7; $r++) {
    $oldLocation = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 6).Value = $oldLocation
}

# New column header
$ws.Range("E1").Value = "Input_group"

# Fill in the new "Input_group" values for each data row
$ws.Range("E2").Value = "2D Media"
$ws.Range("E3").Value = "2D Media"
$ws.Range("E4").Value = "2D Media"
$ws.Range("E5").Value = "Organoid Media"
$ws.Range("E6").Value = "Organoid Media"
$ws.Range("E7").Value = "Organoid Media"

# Widen column B to fit the longer description text
# (49.140625 is the saved OOXML width; the COM layer quantizes ColumnWidth to
# whole pixels, so 48.3 is the input that round-trips closest to that value)
$ws.Columns("B").ColumnWidth = 48.3

# Update the active selection to match the author's final cursor position
$ws.Range("E7").Select()
